# Make elimination language independent (close #406)
#
# The "elimination" metadata column on the Variables sheet's table was
# previously named "en_elimination" and sat at the far right of the table
# (after "en_note"). Make it language independent by renaming it to
# "elimination" and moving it so it sits right after "variable-type"
# (i.e. it becomes the 5th column, pushing "en_variable-label" and
# "en_note" one slot to the right).

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Variables")

# Rename the header of the last table column ("en_elimination" -> "elimination")
# before relocating it.
$ws3.Range("G1").Value = "elimination"

# Move the whole sheet column out of its old slot (G) and drop it back in
# right after "variable-type" (D), i.e. as the new column E. This shifts
# "en_variable-label" (old E) and "en_note" (old F) one column to the right.
$ws3.Columns.Item(7).Cut()
$ws3.Columns.Item(5).Insert()

# Match the selection left behind in the saved file.
$ws3.Range("E2").Select()
